$d = $word.ActiveDocument

# 1. Update the cached TIME field result (the certification date).
$d.Content.Find.Execute("8 de agosto de 2024", $true, $false, $false, $false, $false, $true, 1, $false, "1 de noviembre de 2024", 2) | Out-Null

# 2. Replace the signature block:
#    - "Teniente Coronel JORGE GIOVANNI JIMÉNEZ SÁNCHEZ" / "Subdirector de Meteorología"
#    + "INGRID TATIANA SIERRA GIRALDO" / "Subdirectora de Meteorología" (new paragraphs,
#      new formatting) followed by the old two paragraphs now emptied, and drop the
#      trailing empty paragraph + bookmark that used to wrap the block.
$sigStart = $d.Paragraphs.Item(124).Range.Start
$sigEnd = $d.Paragraphs.Item(128).Range.End
$sigRange = $d.Range($sigStart, $sigEnd)

$newBlockXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Verdana" w:eastAsia="Verdana" w:hAnsi="Verdana" w:cs="Verdana"/><w:b/><w:bCs/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:eastAsia="Verdana" w:hAnsi="Verdana" w:cs="Verdana"/><w:b/><w:bCs/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">INGRID TATIANA SIERRA GIRALDO </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Verdana" w:eastAsia="Verdana" w:hAnsi="Verdana" w:cs="Verdana"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:eastAsia="Verdana" w:hAnsi="Verdana" w:cs="Verdana"/><w:color w:val="000000"/></w:rPr><w:t>Subdirectora de Meteorología</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/><w:contextualSpacing/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="es-CO"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/><w:contextualSpacing/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:lang w:val="es-CO" w:eastAsia="es-CO"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/><w:contextualSpacing/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:lang w:val="es-CO" w:eastAsia="es-CO"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/><w:contextualSpacing/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:lang w:val="es-CO" w:eastAsia="es-CO"/></w:rPr></w:pPr></w:p>
'@

$sigRange.InsertXML($newBlockXml) | Out-Null
